$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update column A (Scanner QR codes) for rows 2..28 ---
$newA = @{
    2  = "https'//scooters.taxify.eu/qr/449-616"
    3  = "https'//scooters.taxify.eu/qr/020-356"
    4  = "https'//scooters.taxify.eu/qr/112-293"
    5  = "https'//scooters.taxify.eu/qr/359-761"
    6  = "https'//scooters.taxify.eu/qr/078-600"
    7  = "https'//scooters.taxify.eu/qr/972-987"
    8  = "https'//scooters.taxify.eu/qr/438-818"
    9  = "https'//scooters.taxify.eu/qr/523-823"
    10 = "https'//scooters.taxify.eu/qr/252-206"
    11 = "https'//scooters.taxify.eu/qr/413-283"
    12 = "https'//scooters.taxify.eu/qr/550-620"
    13 = "https'//scooters.taxify.eu/qr/810-488"
    14 = "https'//scooters.taxify.eu/qr/773-675"
    15 = "https'//scooters.taxify.eu/qr/546-336"
    16 = "https'//scooters.taxify.eu/qr/502-940"
    17 = "https'//scooters.taxify.eu/qr/117-135"
    18 = "https'//scooters.taxify.eu/qr/786-529"
    19 = "https'//scooters.taxify.eu/qr/255-146"
    20 = "https'//scooters.taxify.eu/qr/739-612"
    21 = "https'//scooters.taxify.eu/qr/309-650"
    22 = "https'//scooters.taxify.eu/qr/753-710"
    23 = "https'//scooters.taxify.eu/qr/687-309"
    24 = "https'//scooters.taxify.eu/qr/767-823"
    25 = "https'//scooters.taxify.eu/qr/667-654"
    26 = "https'//scooters.taxify.eu/qr/521-953"
    27 = "https'//scooters.taxify.eu/qr/572-269"
    28 = "https'//scooters.taxify.eu/qr/916-614"
}

foreach ($r in $newA.Keys) {
    $ws.Cells.Item($r, 1).Value = $newA[$r]
}

# --- 2. Extend column B (qr_to_command) formula down through row 217 ---
# Column B already has the formula =RIGHTB(A#,7) for rows 2-6 and 8-102 (as
# shared formulas) plus a stray literal value in B7. Make sure every row from
# 2 to 217 carries the same formula, matching the pre-existing style (s=5)
# used by the surrounding cells.

# First, propagate the existing formatting (style only) from B102 down to the
# newly-populated rows B103:B217 so the new cells keep the right style (s=5)
# instead of picking up a blank/default style.
$ws.Range("B102").Copy()
$ws.Range("B103:B217").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Now (re)apply the formula to every row in the 2..217 range. Existing shared
# formula groups are preserved as-is; new cells get the formula too.
for ($r = 2; $r -le 217; $r++) {
    $ws.Cells.Item($r, 2).Formula = "=RIGHTB(A" + $r + ",7)"
}
